$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "28.307.93"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.560.02"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "210.53"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws.Range("D8") "44.29"
$ws.Range("E8").Value = "  -4.57%  "
Set-TextValue $ws.Range("D9") "23.62"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "1.783.49"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "1.552.40"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "28.294.20"
$ws.Range("E15").Value = "  -0.95%  "
Set-TextValue $ws.Range("D16") "3.64"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  -1.96%  "
Set-TextValue $ws.Range("D18") "60.96"
$ws.Range("E18").Value = "  -1.86%  "
Set-TextValue $ws.Range("D19") "227.69"
$ws.Range("E19").Value = "  -0.43%  "
Set-TextValue $ws.Range("D20") "7.34"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "0.0₃0677"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.18%  "
Set-TextValue $ws.Range("D24") "8.89"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("E25").Value = "  -2.45%  "
Set-TextValue $ws.Range("D26") "150.08"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +2.31%  "
Set-TextValue $ws.Range("D32") "1.06"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "1.377.99"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("E37").Value = "  -3.68%  "
Set-TextValue $ws.Range("D39") "2.63"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("E40").Value = "  -2.17%  "
Set-TextValue $ws.Range("D41") "0.519"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("E43").Value = "  +0.00%  "
Set-TextValue $ws.Range("D44") "0.0470"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  -3.57%  "
Set-TextValue $ws.Range("D47") "62.02"
$ws.Range("E47").Value = "  -1.23%  "
Set-TextValue $ws.Range("D48") "0.916"
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("D49").Value = "1.696.63"
$ws.Range("E49").Value = "  -0.34%  "
Set-TextValue $ws.Range("D50") "85.25"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -2.26%  "
